$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A text cells to stay literal strings ("2025-03-20") rather
# than being auto-parsed into date serial numbers, then drop the
# temporary number-format so no extra cell style is left behind.
$ws.Range("A24:A25").NumberFormat = "@"

$ws.Range("A24").Value = "2025-03-20"
$ws.Range("B24").Value = "substance active"
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 1

$ws.Range("A25").Value = "2025-03-20"
$ws.Range("B25").Value = "éco-régime"
$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 1

$ws.Range("A24:A25").ClearFormats()
